$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update L column (col 12) values
$ws.Range("L7").Value = 15
$ws.Range("L9").Value = 22.5
$ws.Range("L10").Value = 22.5
$ws.Range("L11").Value = 7.5
$ws.Range("L15").Value = 7.5
$ws.Range("L16").Value = 6.3
$ws.Range("L17").Value = 10

# Update M column (col 13) values
$ws.Range("M14").Value = 95
$ws.Range("M18").Value = 95
$ws.Range("M19").Value = 95
$ws.Range("M20").Value = 95
$ws.Range("M21").Value = 95
$ws.Range("M22").Value = 95
$ws.Range("M23").Value = 95
$ws.Range("M24").Value = 95
$ws.Range("M25").Value = 95
$ws.Range("M26").Value = 95
$ws.Range("M27").Value = 95
$ws.Range("M28").Value = 95
$ws.Range("M29").Value = 95
$ws.Range("M30").Value = 95
$ws.Range("M31").Value = 95

# Update N column (col 14) values
$ws.Range("N21").Value = 50
$ws.Range("N22").Value = 50
$ws.Range("N23").Value = 50
$ws.Range("N24").Value = 50
$ws.Range("N25").Value = 50
$ws.Range("N26").Value = 50
$ws.Range("N27").Value = 50
$ws.Range("N28").Value = 50
$ws.Range("N29").Value = 50
$ws.Range("N30").Value = 50
$ws.Range("N31").Value = 50

# Update selection to match new active cell
$ws.Range("W15").Select() | Out-Null
